# "Expert Licences.xlsx" — append a new data row under the header row
# (Account = 12345, Expiry = "2023.12.20 12:23") and leave the new cell
# selected, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 data. Expiry is entered as literal text (not a date/time value).
$ws.Range("A2").Value = 12345
$ws.Range("B2").Value = "2023.12.20 12:23"

# The sheet's row-outline level needs to be bumped to 1 (sheetFormatPr's
# outlineLevelRow) without leaving any row actually grouped. Briefly group a
# spare row below the data and remove it again.
$ws.Range("A3:B3").EntireRow.OutlineLevel = 1
$ws.Range("A3:B3").EntireRow.Delete()

# Final selection lands on the newly entered Expiry cell
$ws.Range("B2").Select()
